$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$ws.Range("L2").Value = 0.8100000000000001
$ws.Range("L3").Value = 1.05
$ws.Range("L4").Value = 1.05
$ws.Range("L5").Value = 0.96
$ws.Range("L6").Value = 0.84
$ws.Range("L7").Value = 0.83
$ws.Range("L8").Value = 0.97
$ws.Range("L9").Value = 1.18
$ws.Range("L10").Value = 0.92
$ws.Range("L11").Value = 1.13
$ws.Range("L12").Value = 1.06
$ws.Range("L13").Value = 0.92
$ws.Range("L14").Value = 0.89
$ws.Range("L15").Value = 1.12
$ws.Range("L16").Value = 1.05
$ws.Range("L17").Value = 1.13
